$d = $word.ActiveDocument

$d.Content.Find.Execute("Фамилия Имя Отчество, БФИ-2202", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Сидорук Данил Вадимович, БФИ-2202", 2)
